$wb = $excel.ActiveWorkbook

$authSheet = $wb.Worksheets.Item("Auth")
$authCaptureSheet = $wb.Worksheets.Item("AuthCapture")

$authSheet.Range("B2").Value = "Thu Oct 27 10:04:04 EDT 2022"
$authSheet.Range("B3").Value = "Thu Oct 27 10:04:27 EDT 2022"
$authSheet.Range("B4").Value = "Thu Oct 27 10:04:43 EDT 2022"
$authSheet.Range("B5").Value = "Thu Oct 27 10:05:01 EDT 2022"
$authSheet.Range("B6").Value = "Thu Oct 27 10:05:17 EDT 2022"

$authCaptureSheet.Range("B2").Value = "Thu Oct 27 10:05:34 EDT 2022"
$authCaptureSheet.Range("B3").Value = "Thu Oct 27 10:05:59 EDT 2022"
$authCaptureSheet.Range("B4").Value = "Thu Oct 27 10:06:23 EDT 2022"
$authCaptureSheet.Range("B5").Value = "Thu Oct 27 10:06:47 EDT 2022"
$authCaptureSheet.Range("B6").Value = "Thu Oct 27 10:07:11 EDT 2022"
